$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shows up in Overview!E/F for both rows and in the "Status" column of
#    both language sheets)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Handback report generated: fill "Latest Target File" (I), "Latest
#    Handback File" (J) and "Latest Handback DateTime" (K) for rows 2 & 3 on
#    both language sheets.
# ---------------------------------------------------------------------------
$targetFileName = "a9ca9bc8-6ed6-41bd-bc60-9c41a806a04e.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aab7feeb224fad3f8e5370c0a9abf18a38e32096/e2e/a9ca9bc8-6ed6-41bd-bc60-9c41a806a04e.md"

$zhCnHandbackFile = "a9ca9bc8-6ed6-41bd-bc60-9c41a806a04e.484e9777cdc6ef738bd963f8dd4327dba3002178.zh-cn.xlf"
$deDeHandbackFile = "a9ca9bc8-6ed6-41bd-bc60-9c41a806a04e.484e9777cdc6ef738bd963f8dd4327dba3002178.de-de.xlf"

$zhCnHandbackTime = "2016-08-28 15:02:44"
$deDeHandbackTime = "2016-08-28 15:02:51"

# zh-cn sheet
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $targetUrl, "", "", $targetFileName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $targetUrl, "", "", $targetFileName)
$wsZhCn.Range("J2").Value = $zhCnHandbackFile
$wsZhCn.Range("J3").Value = $zhCnHandbackFile
$wsZhCn.Range("K2").Value = $zhCnHandbackTime
$wsZhCn.Range("K3").Value = $zhCnHandbackTime

# de-de sheet
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $targetUrl, "", "", $targetFileName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $targetUrl, "", "", $targetFileName)
$wsDeDe.Range("J2").Value = $deDeHandbackFile
$wsDeDe.Range("J3").Value = $deDeHandbackFile
$wsDeDe.Range("K2").Value = $deDeHandbackTime
$wsDeDe.Range("K3").Value = $deDeHandbackTime

# ---------------------------------------------------------------------------
# 3. Widen the columns that now hold longer text (Status on all sheets,
#    Latest Target File / Latest Handback File on the language sheets).
#    ColumnWidth is specified in characters; the stored (raw) width Excel
#    persists is ColumnWidth + 0.8333 (default cell padding), so we dial the
#    requested width in to land on the intended stored widths (~30 and 40).
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.1443713960194   # E -> stored width ~29.98
$wsOverview.Columns.Item(6).ColumnWidth = 29.1443713960194   # F -> stored width ~29.98

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.1443713960194      # C -> stored width ~29.98
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.1666666666667      # I -> stored width 40
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1666666666667      # J -> stored width 40

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.1443713960194      # C -> stored width ~29.98
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.1666666666667      # I -> stored width 40
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1666666666667      # J -> stored width 40
